# Weekly update: insert two new price rows for the most recent week
# (Fecha 44748) at the top of the "Comercializadora del Agro de Limarí -
# Alcachofa" data block, pushing the existing historical rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 198; everything currently at row 198
# onward (the historical data) shifts down to rows 200+.
$ws.Rows("198:199").Insert()

# --- New row 198: Alcachofa, Argentina(o), Primera ---
$ws.Cells.Item(198, 1).Value  = 2
$ws.Cells.Item(198, 2).Value  = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(198, 3).Value  = "Coquimbo"
$ws.Cells.Item(198, 4).Value  = 44748
$ws.Cells.Item(198, 5).Value  = 4
$ws.Cells.Item(198, 6).Value  = 100112013
$ws.Cells.Item(198, 7).Value  = "Alcachofa"
$ws.Cells.Item(198, 8).Value  = "Argentina(o)"
$ws.Cells.Item(198, 9).Value  = "Primera"
$ws.Cells.Item(198, 10).Value = 500
$ws.Cells.Item(198, 11).Value = 11000
$ws.Cells.Item(198, 12).Value = 12000
$ws.Cells.Item(198, 13).Value = 11500
$ws.Cells.Item(198, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(198, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(198, 16).Value = 230
$ws.Cells.Item(198, 17).Value = 50
$ws.Cells.Item(198, 18).Value = "Hortaliza"

# --- New row 199: Alcachofa, Española, Primera ---
$ws.Cells.Item(199, 1).Value  = 2
$ws.Cells.Item(199, 2).Value  = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(199, 3).Value  = "Coquimbo"
$ws.Cells.Item(199, 4).Value  = 44748
$ws.Cells.Item(199, 5).Value  = 4
$ws.Cells.Item(199, 6).Value  = 100112013
$ws.Cells.Item(199, 7).Value  = "Alcachofa"
$ws.Cells.Item(199, 8).Value  = "Española"
$ws.Cells.Item(199, 9).Value  = "Primera"
$ws.Cells.Item(199, 10).Value = 800
$ws.Cells.Item(199, 11).Value = 14000
$ws.Cells.Item(199, 12).Value = 15000
$ws.Cells.Item(199, 13).Value = 14500
$ws.Cells.Item(199, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(199, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(199, 16).Value = 483
$ws.Cells.Item(199, 17).Value = 30
$ws.Cells.Item(199, 18).Value = "Hortaliza"

# Make sure the date cells keep the workbook's date number format
# (same as every other "Fecha" cell in column D).
$ws.Range("D198:D199").NumberFormat = $ws.Range("D200").NumberFormat

$rowCount = $ws.UsedRange.Rows.Count
"Done: inserted rows 198-199, used range now has $rowCount rows"
